$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.804.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.888.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7657'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.99'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3116'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.24'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -7.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07203'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08068'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7635'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.503'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.904.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.17'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.133'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.820.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007750'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.157.91'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.102'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1551'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.381'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.07'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.71'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.035'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.438'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.549'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.455'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.094'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05477'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.256'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7451'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01915'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.780'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.142.71'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +11.21%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4406'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '73.27'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.880'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8497'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.68'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.883'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.955'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.436'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.97%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.015'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +10.31%  '
